# Apply the "multiple channels and transformations" update to the
# experiment workbook (beads / cells sheets).

$wb = $excel.ActiveWorkbook
$beads = $wb.Worksheets.Item("beads")
$cells = $wb.Worksheets.Item("cells")

# --- beads sheet: rename the MEFL column header -------------------------
$beads.Range("B1").Value = "FL1-H Peaks"

# --- cells sheet: insert a new "FL1-H Transform" column after B --------
$cells.Columns.Item(3).Insert()
$cells.Columns.Item(3).ColumnWidth = $cells.Columns.Item(1).ColumnWidth

$cells.Range("C1").Value = "FL1-H Transform"
$cells.Range("C2").Clear()
$cells.Range("C3").Value = "None"
$cells.Range("C4").Value = "Exponential"
$cells.Range("C5").Value = "Mef"
$cells.Range("C6").Value = "Mef"

# --- view state: "cells" becomes the active/selected tab ---------------
$beads.Range("B2").Select()

$cells.Activate()
$cells.Range("E28").Select()
